# Generate Report for Archive
$wb = $excel.ActiveWorkbook

# Update the "Ready for handoff" status text to "In Translation" wherever it
# appears (Overview sheet columns E/F, and the per-language sheets' Status column).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $targets = @()
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Text) {
            $targets += $cell.Address()
        }
    }
    foreach ($addr in $targets) {
        $ws.Range($addr).Value = "In Translation"
    }
}

# Narrow the "Status" columns that previously held the "Ready for handoff" text.
# (Target stored column width is 13.4101845877511 "characters" of display
# width; the host's ColumnWidth setter only has 1/6-character resolution, so
# 12.5 is the input that lands on the closest representable stored width.)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
